$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(2, 0.0239616380717946, 0.06344570286477988),
    @(3, -0.008678078645303806, -0.007630201523759542),
    @(4, -0.001854580733088329, 0.02057698068935901),
    @(5, -0.05225139074879027, -0.1902565838498182),
    @(6, 0.01537564809319974, 0.1530303873050321),
    @(7, 0.3466348215901034, 0.261629465054335),
    @(8, 0.278357800202353, 0.2585824801839264),
    @(9, 0.4225003425708806, -0.1239312759104667),
    @(10, 0.7061988505389594, -0.03952148328333436),
    @(11, -0.01123539325197448, 0.05157341663780133),
    @(12, 0.00484950103050021, 0.3289914647728868),
    @(13, 0.04582151089640971, -0.0269968212924624),
    @(14, 0.212782049431686, 0.08714304692436992),
    @(15, -0.139487960292954, 0.7222264592375381),
    @(16, 0.2091404642647843, 0.1130721382794221),
    @(17, -0.1039045934638865, 0.360592639479706)
)

foreach ($row in $values) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
